# Apply the "input system override path" key/value rows on Sheet1:
# B5 and B6 move from raw numeric placeholders to real string values
# (Korean greeting text), and the active selection shifts to C7.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("B5").Value = "안녕안녕안녕안녕안녕안녕안녕안녕안녕안녕안녕안녕안녕안녕안녕안녕안녕안녕안녕안녕안녕안녕안녕안녕안녕안녕안녕안녕안녕안녕안녕안녕안녕안녕"
$ws.Range("B6").Value = "반가워반가워반가워반가워반가워반가워반가워반가워반가워반가워반가워반가워반가워반가워반가워반가워반가워반가워반가워반가워반가워반가워반가워반가워"

$ws.Range("C7").Select()
